$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> list of (row, col, value) updates derived from the
# upstream Sheets/Mandragora_Profits.xlsx diff (refreshed market-board pricing data).

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 967.8570999999999  # H17: was 1004.13336
$ws.Cells.Item(17, 10).Value = 985.1852  # J17: was 1021.5172
$ws.Cells.Item(17, 12).Value = 2955.5556  # L17: was 3064.5516
$ws.Cells.Item(17, 14).Value = -3291.5556  # N17: was -3400.5516

# Row 93
$ws.Cells.Item(93, 8).Value = 49999.5  # H93: was 54800
$ws.Cells.Item(93, 10).Value = 49999.5  # J93: was 54800
$ws.Cells.Item(93, 12).Value = 49999.5  # L93: was 54800
$ws.Cells.Item(93, 14).Value = -54991.5  # N93: was -59792

# Row 116
$ws.Cells.Item(116, 8).Value = 3714.1428  # H116: was 3341.6365
$ws.Cells.Item(116, 9).Value = 3000  # I116: was 2959.6
$ws.Cells.Item(116, 10).Value = 3833.1667  # J116: was 3660
$ws.Cells.Item(116, 11).Value = 3000  # K116: was 2959.6
$ws.Cells.Item(116, 12).Value = 3833.1667  # L116: was 3660
$ws.Cells.Item(116, 13).Value = 442  # M116: was 482.4000000000001
$ws.Cells.Item(116, 14).Value = -10717.1667  # N116: was -10544

# Row 127
$ws.Cells.Item(127, 8).Value = 1292.875  # H127: was 1304.3
$ws.Cells.Item(127, 9).Value = 1218.6  # I127: was 1182.1666
$ws.Cells.Item(127, 10).Value = 1416.6666  # J127: was 1487.5
$ws.Cells.Item(127, 11).Value = 3655.8  # K127: was 3546.4998
$ws.Cells.Item(127, 12).Value = 4249.9998  # L127: was 4462.5
$ws.Cells.Item(127, 13).Value = 1304.2  # M127: was 1413.5002
$ws.Cells.Item(127, 14).Value = -14169.9998  # N127: was -14382.5

# Row 129
$ws.Cells.Item(129, 8).Value = 1642.6052  # H129: was 1703.8206
$ws.Cells.Item(129, 9).Value = 690.9167  # I129: was 717.36365
$ws.Cells.Item(129, 10).Value = 2081.8462  # J129: was 2091.3572
$ws.Cells.Item(129, 11).Value = 2072.7501  # K129: was 2152.09095
$ws.Cells.Item(129, 12).Value = 6245.5386  # L129: was 6274.071599999999
$ws.Cells.Item(129, 13).Value = 2927.2499  # M129: was 2847.90905
$ws.Cells.Item(129, 14).Value = -16245.5386  # N129: was -16274.0716

# Row 137
$ws.Cells.Item(137, 8).Value = 1609.2826  # H137: was 1382.4237
$ws.Cells.Item(137, 9).Value = 3138.1538  # I137: was 1955.5416
$ws.Cells.Item(137, 10).Value = 1007  # J137: was 989.4286
$ws.Cells.Item(137, 11).Value = 9414.4614  # K137: was 5866.6248
$ws.Cells.Item(137, 12).Value = 3021  # L137: was 2968.2858
$ws.Cells.Item(137, 13).Value = -6864.4614  # M137: was -3316.6248
$ws.Cells.Item(137, 14).Value = -8121  # N137: was -8068.2858

# Row 138
$ws.Cells.Item(138, 8).Value = 1525.3684  # H138: was 1600.5596
$ws.Cells.Item(138, 9).Value = 1150.7925  # I138: was 1147.7963
$ws.Cells.Item(138, 10).Value = 2388.5217  # J138: was 2415.5334
$ws.Cells.Item(138, 11).Value = 3452.3775  # K138: was 3443.3889
$ws.Cells.Item(138, 12).Value = 7165.5651  # L138: was 7246.600199999999
$ws.Cells.Item(138, 13).Value = 1687.6225  # M138: was 1696.6111
$ws.Cells.Item(138, 14).Value = -17445.5651  # N138: was -17526.6002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6904.273  # H32: was 5808.81
$ws.Cells.Item(32, 9).Value = 5831.728  # I32: was 4505.066
$ws.Cells.Item(32, 10).Value = 21000.572  # J32: was 18991.111
$ws.Cells.Item(32, 11).Value = 5831.728  # K32: was 4505.066
$ws.Cells.Item(32, 12).Value = 21000.572  # L32: was 18991.111
$ws.Cells.Item(32, 13).Value = -5544.728  # M32: was -4218.066
$ws.Cells.Item(32, 14).Value = -21574.572  # N32: was -19565.111

# Row 109
$ws.Cells.Item(109, 8).Value = 14951.5  # H109: was 15590
$ws.Cells.Item(109, 10).Value = 14951.5  # J109: was 15590
$ws.Cells.Item(109, 12).Value = 14951.5  # L109: was 15590
$ws.Cells.Item(109, 14).Value = -17725.5  # N109: was -18364

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1490.2526  # H31: was 1818.39
$ws.Cells.Item(31, 9).Value = 991.597  # I31: was 1106.1552
$ws.Cells.Item(31, 10).Value = 2534.3125  # J31: was 2801.9524
$ws.Cells.Item(31, 11).Value = 991.597  # K31: was 1106.1552
$ws.Cells.Item(31, 12).Value = 2534.3125  # L31: was 2801.9524
$ws.Cells.Item(31, 13).Value = -696.597  # M31: was -811.1551999999999
$ws.Cells.Item(31, 14).Value = -3124.3125  # N31: was -3391.9524

# Row 34
$ws.Cells.Item(34, 8).Value = 1490.2526  # H34: was 1818.39
$ws.Cells.Item(34, 9).Value = 991.597  # I34: was 1106.1552
$ws.Cells.Item(34, 10).Value = 2534.3125  # J34: was 2801.9524
$ws.Cells.Item(34, 11).Value = 991.597  # K34: was 1106.1552
$ws.Cells.Item(34, 12).Value = 2534.3125  # L34: was 2801.9524
$ws.Cells.Item(34, 13).Value = -789.597  # M34: was -904.1551999999999
$ws.Cells.Item(34, 14).Value = -2938.3125  # N34: was -3205.9524

# Row 36
$ws.Cells.Item(36, 8).Value = 0  # H36: was 1663
$ws.Cells.Item(36, 9).Value = 0  # I36: was 1663
$ws.Cells.Item(36, 11).Value = 0  # K36: was 1663
$ws.Cells.Item(36, 13).ClearContents()  # M36: was -1275

# Row 40
$ws.Cells.Item(40, 8).Value = 0  # H40: was 1663
$ws.Cells.Item(40, 9).Value = 0  # I40: was 1663
$ws.Cells.Item(40, 11).Value = 0  # K40: was 1663
$ws.Cells.Item(40, 13).ClearContents()  # M40: was -1503

# Row 74
$ws.Cells.Item(74, 8).Value = 22657  # H74: was 22683
$ws.Cells.Item(74, 10).Value = 22657  # J74: was 22683
$ws.Cells.Item(74, 12).Value = 22657  # L74: was 22683
$ws.Cells.Item(74, 14).Value = -24405  # N74: was -24431

# Row 77
$ws.Cells.Item(77, 8).Value = 22657  # H77: was 22683
$ws.Cells.Item(77, 10).Value = 22657  # J77: was 22683
$ws.Cells.Item(77, 12).Value = 67971  # L77: was 68049
$ws.Cells.Item(77, 14).Value = -76707  # N77: was -76785

# Row 88
$ws.Cells.Item(88, 8).Value = 25387.375  # H88: was 23624.875
$ws.Cells.Item(88, 9).Value = 12400  # I88: was 12500
$ws.Cells.Item(88, 10).Value = 29716.5  # J88: was 27333.166
$ws.Cells.Item(88, 11).Value = 12400  # K88: was 12500
$ws.Cells.Item(88, 12).Value = 29716.5  # L88: was 27333.166
$ws.Cells.Item(88, 13).Value = -11994  # M88: was -12094
$ws.Cells.Item(88, 14).Value = -30528.5  # N88: was -28145.166

# Row 91
$ws.Cells.Item(91, 8).Value = 25387.375  # H91: was 23624.875
$ws.Cells.Item(91, 9).Value = 12400  # I91: was 12500
$ws.Cells.Item(91, 10).Value = 29716.5  # J91: was 27333.166
$ws.Cells.Item(91, 11).Value = 12400  # K91: was 12500
$ws.Cells.Item(91, 12).Value = 29716.5  # L91: was 27333.166
$ws.Cells.Item(91, 13).Value = -10996  # M91: was -11096
$ws.Cells.Item(91, 14).Value = -32524.5  # N91: was -30141.166

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Cells.Item(23, 8).Value = 176.07143  # H23: was 182.6923
$ws.Cells.Item(23, 10).Value = 209.3  # J23: was 222.55556
$ws.Cells.Item(23, 12).Value = 627.9000000000001  # L23: was 667.66668
$ws.Cells.Item(23, 14).Value = -1097.9  # N23: was -1137.66668

# Row 68
$ws.Cells.Item(68, 8).Value = 1239.3636  # H68: was 1011.5
$ws.Cells.Item(68, 9).Value = 883.3333  # I68: was 870.2857
$ws.Cells.Item(68, 10).Value = 1372.875  # J68: was 2000
$ws.Cells.Item(68, 11).Value = 2649.9999  # K68: was 2610.8571
$ws.Cells.Item(68, 12).Value = 4118.625  # L68: was 6000
$ws.Cells.Item(68, 13).Value = -1838.9999  # M68: was -1799.8571
$ws.Cells.Item(68, 14).Value = -5740.625  # N68: was -7622

# Row 71
$ws.Cells.Item(71, 8).Value = 1239.3636  # H71: was 1011.5
$ws.Cells.Item(71, 9).Value = 883.3333  # I71: was 870.2857
$ws.Cells.Item(71, 10).Value = 1372.875  # J71: was 2000
$ws.Cells.Item(71, 11).Value = 7949.9997  # K71: was 7832.571300000001
$ws.Cells.Item(71, 12).Value = 12355.875  # L71: was 18000
$ws.Cells.Item(71, 13).Value = -3893.9997  # M71: was -3776.571300000001
$ws.Cells.Item(71, 14).Value = -20467.875  # N71: was -26112

# Row 132
$ws.Cells.Item(132, 8).Value = 1331.4615  # H132: was 1175.25
$ws.Cells.Item(132, 9).Value = 782.3333  # I132: was 922
$ws.Cells.Item(132, 10).Value = 1802.1428  # J132: was 1935
$ws.Cells.Item(132, 11).Value = 7040.9997  # K132: was 8298
$ws.Cells.Item(132, 12).Value = 16219.2852  # L132: was 17415
$ws.Cells.Item(132, 13).Value = -4510.9997  # M132: was -5768
$ws.Cells.Item(132, 14).Value = -21279.2852  # N132: was -22475

$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Cells.Item(64, 8).Value = 8136  # H64: was 15068
$ws.Cells.Item(64, 10).Value = 0  # J64: was 22000
$ws.Cells.Item(64, 12).Value = 0  # L64: was 22000
$ws.Cells.Item(64, 14).ClearContents()  # N64: was -22450

# Row 67
$ws.Cells.Item(67, 8).Value = 8136  # H67: was 15068
$ws.Cells.Item(67, 10).Value = 0  # J67: was 22000
$ws.Cells.Item(67, 12).Value = 0  # L67: was 22000
$ws.Cells.Item(67, 14).ClearContents()  # N67: was -23560

# Row 87
$ws.Cells.Item(87, 8).Value = 14250  # H87: was 0
$ws.Cells.Item(87, 9).Value = 14500  # I87: was 0
$ws.Cells.Item(87, 10).Value = 14000  # J87: was 0
$ws.Cells.Item(87, 11).Value = 14500  # K87: was 0
$ws.Cells.Item(87, 12).Value = 14000  # L87: was 0
$ws.Cells.Item(87, 13).Value = -13377  # M87: was blank
$ws.Cells.Item(87, 14).Value = -16246  # N87: was blank

# Row 88
$ws.Cells.Item(88, 8).Value = 13230  # H88: was 9500
$ws.Cells.Item(88, 9).Value = 7800  # I88: was 1000
$ws.Cells.Item(88, 10).Value = 15945  # J88: was 18000
$ws.Cells.Item(88, 11).Value = 7800  # K88: was 1000
$ws.Cells.Item(88, 12).Value = 15945  # L88: was 18000
$ws.Cells.Item(88, 13).Value = -7372  # M88: was -572
$ws.Cells.Item(88, 14).Value = -16801  # N88: was -18856

# Row 90
$ws.Cells.Item(90, 8).Value = 14250  # H90: was 0
$ws.Cells.Item(90, 9).Value = 14500  # I90: was 0
$ws.Cells.Item(90, 10).Value = 14000  # J90: was 0
$ws.Cells.Item(90, 11).Value = 43500  # K90: was 0
$ws.Cells.Item(90, 12).Value = 42000  # L90: was 0
$ws.Cells.Item(90, 13).Value = -37884  # M90: was blank
$ws.Cells.Item(90, 14).Value = -53232  # N90: was blank

# Row 91
$ws.Cells.Item(91, 8).Value = 13230  # H91: was 9500
$ws.Cells.Item(91, 9).Value = 7800  # I91: was 1000
$ws.Cells.Item(91, 10).Value = 15945  # J91: was 18000
$ws.Cells.Item(91, 11).Value = 7800  # K91: was 1000
$ws.Cells.Item(91, 12).Value = 15945  # L91: was 18000
$ws.Cells.Item(91, 13).Value = -6318  # M91: was 482
$ws.Cells.Item(91, 14).Value = -18909  # N91: was -20964

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Cells.Item(63, 8).Value = 24083.166  # H63: was 21120.834
$ws.Cells.Item(63, 9).Value = 0  # I63: was 2226
$ws.Cells.Item(63, 10).Value = 24083.166  # J63: was 24899.8
$ws.Cells.Item(63, 11).Value = 0  # K63: was 2226
$ws.Cells.Item(63, 12).Value = 24083.166  # L63: was 24899.8
$ws.Cells.Item(63, 13).ClearContents()  # M63: was -1602
$ws.Cells.Item(63, 14).Value = -25331.166  # N63: was -26147.8

# Row 66
$ws.Cells.Item(66, 8).Value = 24083.166  # H66: was 21120.834
$ws.Cells.Item(66, 9).Value = 0  # I66: was 2226
$ws.Cells.Item(66, 10).Value = 24083.166  # J66: was 24899.8
$ws.Cells.Item(66, 11).Value = 0  # K66: was 6678
$ws.Cells.Item(66, 12).Value = 72249.49800000001  # L66: was 74699.39999999999
$ws.Cells.Item(66, 13).ClearContents()  # M66: was -3558
$ws.Cells.Item(66, 14).Value = -78489.49800000001  # N66: was -80939.39999999999

# Row 75
$ws.Cells.Item(75, 8).Value = 15000  # H75: was 14450
$ws.Cells.Item(75, 9).Value = 15000  # I75: was 0
$ws.Cells.Item(75, 10).Value = 15000  # J75: was 14450
$ws.Cells.Item(75, 11).Value = 15000  # K75: was 0
$ws.Cells.Item(75, 12).Value = 15000  # L75: was 14450
$ws.Cells.Item(75, 13).Value = -14064  # M75: was blank
$ws.Cells.Item(75, 14).Value = -16872  # N75: was -16322

# Row 78
$ws.Cells.Item(78, 8).Value = 15000  # H78: was 14450
$ws.Cells.Item(78, 9).Value = 15000  # I78: was 0
$ws.Cells.Item(78, 10).Value = 15000  # J78: was 14450
$ws.Cells.Item(78, 11).Value = 45000  # K78: was 0
$ws.Cells.Item(78, 12).Value = 45000  # L78: was 43350
$ws.Cells.Item(78, 13).Value = -40320  # M78: was blank
$ws.Cells.Item(78, 14).Value = -54360  # N78: was -52710

# Row 82
$ws.Cells.Item(82, 8).Value = 9000  # H82: was 0
$ws.Cells.Item(82, 10).Value = 9000  # J82: was 0
$ws.Cells.Item(82, 12).Value = 9000  # L82: was 0
$ws.Cells.Item(82, 14).Value = -9766  # N82: was blank

# Row 85
$ws.Cells.Item(85, 8).Value = 9000  # H85: was 0
$ws.Cells.Item(85, 10).Value = 9000  # J85: was 0
$ws.Cells.Item(85, 12).Value = 9000  # L85: was 0
$ws.Cells.Item(85, 14).Value = -11652  # N85: was blank

# Row 86
$ws.Cells.Item(86, 8).Value = 14865  # H86: was 17081.25
$ws.Cells.Item(86, 10).Value = 14865  # J86: was 17081.25
$ws.Cells.Item(86, 12).Value = 14865  # L86: was 17081.25
$ws.Cells.Item(86, 14).Value = -17111  # N86: was -19327.25

# Row 89
$ws.Cells.Item(89, 8).Value = 14865  # H89: was 17081.25
$ws.Cells.Item(89, 10).Value = 14865  # J89: was 17081.25
$ws.Cells.Item(89, 12).Value = 74325  # L89: was 85406.25
$ws.Cells.Item(89, 14).Value = -85557  # N89: was -96638.25

# Row 92
$ws.Cells.Item(92, 8).Value = 29000  # H92: was 30549
$ws.Cells.Item(92, 10).Value = 29000  # J92: was 30549
$ws.Cells.Item(92, 12).Value = 29000  # L92: was 30549
$ws.Cells.Item(92, 14).Value = -33992  # N92: was -35541

